$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.449.52'
$ws.Range('E2').Value = '  -1.52%  '
$ws.Range('D3').Value = '2.184.88'
$ws.Range('E3').Value = '  -2.30%  '
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '251.99'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.34%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.606'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.07%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '74.71'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.34%  '
$ws.Range('E8').Value = '  -0.08%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.581'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -5.41%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '40.28'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.60%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0910'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.62%  '
$ws.Range('E12').Value = '  -0.75%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.76'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.95%  '
$ws.Range('D14').Value = '2.513.26'
$ws.Range('E14').Value = '  -2.38%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.20'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.88%  '
$ws.Range('D16').Value = '2.174.12'
$ws.Range('E16').Value = '  -2.56%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.771'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -5.47%  '
$ws.Range('D18').Value = '42.371.61'
$ws.Range('E18').Value = '  -1.50%  '
$ws.Range('E19').Value = '  -3.04%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '70.81'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.23%  '
$ws.Range('E21').Value = '  -2.21%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '226.88'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.63%  '
$ws.Range('E23').Value = '  -10.28%  '
$ws.Range('E24').Value = '  -3.20%  '
$ws.Range('E25').Value = '  +0.02%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '10.45'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -4.74%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '3.39'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.11%  '
$ws.Range('E28').Value = '  -0.03%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.16'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -4.01%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '172.25'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.50%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '37.01'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.92%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '20.02'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.61%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0815'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.48%  '
$ws.Range('E34').Value = '  -4.31%  '
$ws.Range('E35').Value = '  -1.90%  '
$ws.Range('E36').Value = '  -4.56%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.20'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.22%  '
$ws.Range('E38').Value = '  -0.40%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '12.05'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -8.70%  '
$ws.Range('E40').Value = '  -3.34%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.58'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +11.52%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.17'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -7.71%  '
$ws.Range('B43').Value = 'MultiversX'
$ws.Range('C43').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '58.63'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.66%  '
$ws.Range('B44').Value = 'Algorand'
$ws.Range('C44').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.193'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -3.16%  '
$ws.Range('E45').Value = '  -3.65%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0970'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.55%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.17'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -4.05%  '
$ws.Range('E48').Value = '  +3.09%  '
$ws.Range('E49').Value = '  -2.03%  '
$ws.Range('E50').Value = '  -2.73%  '
$ws.Range('E51').Value = '  -0.74%  '
